# Weekly fruit/vegetable price update:
# A new record (week) is inserted as row 48, shifting all subsequent
# existing records (old rows 48-74) down by one row to rows 49-75.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(48).Insert()

$ws.Cells.Item(48, 1).Value  = 10
$ws.Cells.Item(48, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(48, 3).Value  = "La Araucanía"
$ws.Cells.Item(48, 4).Value  = 44455
$ws.Cells.Item(48, 5).Value  = 9
$ws.Cells.Item(48, 6).Value  = 100112031
$ws.Cells.Item(48, 7).Value  = "Poroto verde"
$ws.Cells.Item(48, 8).Value  = "Sin especificar"
$ws.Cells.Item(48, 9).Value  = "Primera"
$ws.Cells.Item(48, 10).Value = 20
$ws.Cells.Item(48, 11).Value = 40000
$ws.Cells.Item(48, 12).Value = 40000
$ws.Cells.Item(48, 13).Value = 40000
$ws.Cells.Item(48, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(48, 15).Value = "Perú"
$ws.Cells.Item(48, 16).Value = 1600
$ws.Cells.Item(48, 17).Value = 25
$ws.Cells.Item(48, 18).Value = "Hortaliza"
